# Fruta / hortaliza, semanal
# Insert 3 new weekly rows at the top of the Femacal de La Calera - Kiwi
# data block (rows 958:960), pushing the existing data down by 3 rows.
# The three new rows mirror the row that lands directly below them
# (same date/quality/price-unit), except for a refreshed Volumen (col M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 958:997 down by inserting 3 blank rows at row 958.
$ws.Range("A958:T960").EntireRow.Insert()

# New row 958 (mirrors the row now at 961, Volumen 56 -> 75)
$ws.Range("A958").Value = 3
$ws.Range("B958").Value = "Femacal de La Calera"
$ws.Range("C958").Value = "Coquimbo"
$ws.Range("D958").Value = 44705
$ws.Range("E958").Value = 5
$ws.Range("F958").Value = "Fruta"
$ws.Range("G958").Value = 100101
$ws.Range("H958").Value = "Berries"
$ws.Range("I958").Value = 100101007
$ws.Range("J958").Value = "Kiwi"
$ws.Range("K958").Value = "Hayward"
$ws.Range("L958").Value = "Especial"
$ws.Range("M958").Value = 75
$ws.Range("N958").Value = 8000
$ws.Range("O958").Value = 8000
$ws.Range("P958").Value = 8000
$ws.Range("Q958").Value = "$/bandeja 10 kilos"
$ws.Range("R958").Value = "Región de O'Higgins"
$ws.Range("S958").Value = 800
$ws.Range("T958").Value = 10

# New row 959 (mirrors the row now at 962, Volumen 67 -> 78)
$ws.Range("A959").Value = 3
$ws.Range("B959").Value = "Femacal de La Calera"
$ws.Range("C959").Value = "Coquimbo"
$ws.Range("D959").Value = 44705
$ws.Range("E959").Value = 5
$ws.Range("F959").Value = "Fruta"
$ws.Range("G959").Value = 100101
$ws.Range("H959").Value = "Berries"
$ws.Range("I959").Value = 100101007
$ws.Range("J959").Value = "Kiwi"
$ws.Range("K959").Value = "Hayward"
$ws.Range("L959").Value = "Primera"
$ws.Range("M959").Value = 78
$ws.Range("N959").Value = 7000
$ws.Range("O959").Value = 7000
$ws.Range("P959").Value = 7000
$ws.Range("Q959").Value = "$/bandeja 10 kilos"
$ws.Range("R959").Value = "Región de O'Higgins"
$ws.Range("S959").Value = 700
$ws.Range("T959").Value = 10

# New row 960 (mirrors the row now at 963, Volumen 60 -> 70)
$ws.Range("A960").Value = 3
$ws.Range("B960").Value = "Femacal de La Calera"
$ws.Range("C960").Value = "Coquimbo"
$ws.Range("D960").Value = 44705
$ws.Range("E960").Value = 5
$ws.Range("F960").Value = "Fruta"
$ws.Range("G960").Value = 100101
$ws.Range("H960").Value = "Berries"
$ws.Range("I960").Value = 100101007
$ws.Range("J960").Value = "Kiwi"
$ws.Range("K960").Value = "Hayward"
$ws.Range("L960").Value = "Segunda"
$ws.Range("M960").Value = 70
$ws.Range("N960").Value = 6000
$ws.Range("O960").Value = 6000
$ws.Range("P960").Value = 6000
$ws.Range("Q960").Value = "$/bandeja 10 kilos"
$ws.Range("R960").Value = "Región de O'Higgins"
$ws.Range("S960").Value = 600
$ws.Range("T960").Value = 10
